$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 54203
$ws.Range("J7").Value = 54203
$ws.Range("L7").Value = 54203
$ws.Range("N7").Value = -54427

$ws.Range("H14").Value = 54203
$ws.Range("J14").Value = 54203
$ws.Range("L14").Value = 54203
$ws.Range("N14").Value = -54585

$ws.Range("H40").Value = 3362.5
$ws.Range("J40").Value = 3816.6667
$ws.Range("L40").Value = 3816.6667
$ws.Range("N40").Value = -4166.6667

$ws.Range("H64").Value = 3350
$ws.Range("J64").Value = 3300
$ws.Range("L64").Value = 3300
$ws.Range("N64").Value = -3796

$ws.Range("H67").Value = 3350
$ws.Range("J67").Value = 3300
$ws.Range("L67").Value = 3300
$ws.Range("N67").Value = -5016

$ws.Range("H74").Value = 4345.6665
$ws.Range("I74").Value = 3575
$ws.Range("J74").Value = 5887
$ws.Range("K74").Value = 3575
$ws.Range("L74").Value = 5887
$ws.Range("M74").Value = -2639
$ws.Range("N74").Value = -7759

$ws.Range("H77").Value = 4345.6665
$ws.Range("I77").Value = 3575
$ws.Range("J77").Value = 5887
$ws.Range("K77").Value = 17875
$ws.Range("L77").Value = 29435
$ws.Range("M77").Value = -13195
$ws.Range("N77").Value = -38795

$ws.Range("H116").Value = 3181.3333
$ws.Range("I116").Value = 2900
$ws.Range("J116").Value = 3462.6667
$ws.Range("K116").Value = 2900
$ws.Range("L116").Value = 3462.6667
$ws.Range("M116").Value = 542
$ws.Range("N116").Value = -10346.6667

$ws.Range("H138").Value = 2203.8928
$ws.Range("I138").Value = 2612.25
$ws.Range("J138").Value = 2107.8088
$ws.Range("K138").Value = 7836.75
$ws.Range("L138").Value = 6323.426399999999
$ws.Range("M138").Value = -2696.75
$ws.Range("N138").Value = -16603.4264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1670.0358
$ws.Range("I74").Value = 1316.6
$ws.Range("J74").Value = 2077.8462
$ws.Range("K74").Value = 1316.6
$ws.Range("L74").Value = 2077.8462
$ws.Range("M74").Value = -442.5999999999999
$ws.Range("N74").Value = -3825.8462

$ws.Range("H77").Value = 1670.0358
$ws.Range("I77").Value = 1316.6
$ws.Range("J77").Value = 2077.8462
$ws.Range("K77").Value = 6583
$ws.Range("L77").Value = 10389.231
$ws.Range("M77").Value = -2215
$ws.Range("N77").Value = -19125.231

$ws.Range("H132").Value = 3170.0193
$ws.Range("I132").Value = 2685.457
$ws.Range("J132").Value = 4167.647
$ws.Range("K132").Value = 8056.370999999999
$ws.Range("L132").Value = 12502.941
$ws.Range("M132").Value = -5526.370999999999
$ws.Range("N132").Value = -17562.941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2757.3447
$ws.Range("I134").Value = 2772.1667
$ws.Range("J134").Value = 2733.0908
$ws.Range("K134").Value = 8316.500100000001
$ws.Range("L134").Value = 8199.2724
$ws.Range("M134").Value = -5781.500100000001
$ws.Range("N134").Value = -13269.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50588

$ws.Range("H132").Value = 11113128
$ws.Range("I132").Value = 1324.3
$ws.Range("K132").Value = 3972.9
$ws.Range("M132").Value = -1442.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1710.0294
$ws.Range("I16").Value = 766.3333
$ws.Range("J16").Value = 1801.3549
$ws.Range("K16").Value = 2298.9999
$ws.Range("L16").Value = 5404.0647
$ws.Range("M16").Value = -2125.9999
$ws.Range("N16").Value = -5750.0647

$ws.Range("H110").Value = 11501.345
$ws.Range("J110").Value = 12327.385
$ws.Range("L110").Value = 36982.155
$ws.Range("N110").Value = -45162.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 255.71428
$ws.Range("I107").Value = 220
$ws.Range("K107").Value = 220
$ws.Range("M107").Value = 1700

$ws.Range("H132").Value = 2262.3242
$ws.Range("I132").Value = 1827.4642
$ws.Range("K132").Value = 5482.392599999999
$ws.Range("M132").Value = -2952.392599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1002
$ws.Range("I16").Value = 771.61536
$ws.Range("J16").Value = 2499.5
$ws.Range("K16").Value = 771.61536
$ws.Range("L16").Value = 2499.5
$ws.Range("M16").Value = -601.61536
$ws.Range("N16").Value = -2839.5

$ws.Range("H68").Value = 2064.145
$ws.Range("I68").Value = 1859.5883
$ws.Range("J68").Value = 2262.8572
$ws.Range("K68").Value = 1859.5883
$ws.Range("L68").Value = 2262.8572
$ws.Range("M68").Value = -1110.5883
$ws.Range("N68").Value = -3760.8572

$ws.Range("H71").Value = 2064.145
$ws.Range("I71").Value = 1859.5883
$ws.Range("J71").Value = 2262.8572
$ws.Range("K71").Value = 9297.941499999999
$ws.Range("L71").Value = 11314.286
$ws.Range("M71").Value = -5553.941499999999
$ws.Range("N71").Value = -18802.286

$ws.Range("H93").Value = 10253.75
$ws.Range("I93").Value = 12338.333
$ws.Range("K93").Value = 12338.333
$ws.Range("M93").Value = -11090.333

$ws.Range("H132").Value = 4934.25
$ws.Range("I132").Value = 4174.143
$ws.Range("K132").Value = 12522.429
$ws.Range("M132").Value = -9992.429

$ws.Range("H136").Value = 3877287.2
$ws.Range("I136").Value = 963.6429000000001
$ws.Range("J136").Value = 11113091
$ws.Range("K136").Value = 2890.9287
$ws.Range("L136").Value = 33339273
$ws.Range("M136").Value = -340.9287000000004
$ws.Range("N136").Value = -33344373

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5569.909
$ws.Range("I81").Value = 9690
$ws.Range("J81").Value = 4024.875
$ws.Range("K81").Value = 19380
$ws.Range("L81").Value = 8049.75
$ws.Range("M81").Value = -18319
$ws.Range("N81").Value = -10171.75

$ws.Range("H84").Value = 5569.909
$ws.Range("I84").Value = 9690
$ws.Range("J84").Value = 4024.875
$ws.Range("K84").Value = 96900
$ws.Range("L84").Value = 40248.75
$ws.Range("M84").Value = -91596
$ws.Range("N84").Value = -50856.75

$ws.Range("H107").Value = 668.3333
$ws.Range("I107").Value = 751
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 2253
$ws.Range("L107").Value = 1509
$ws.Range("M107").Value = -333
$ws.Range("N107").Value = -5349

$ws.Range("H136").Value = 2885.0688
$ws.Range("I136").Value = 2687.4736
$ws.Range("J136").Value = 3260.5
$ws.Range("K136").Value = 8062.4208
$ws.Range("L136").Value = 9781.5
$ws.Range("M136").Value = -5512.4208
$ws.Range("N136").Value = -14881.5
